# Weekly update: a new "Fruta, Macroferia Regional de Talca - Arándano (blue)"
# record is published, so a new row is inserted just above the most recent
# existing record (row 78) and the sheet's older rows shift down by one.
# The brand-new row is seeded with the same data as the record that is now
# immediately below it, except it carries the new reporting date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 78

# Push the existing data (old rows 78..102) down to rows 79..103.
$ws.Rows($newRow).Insert()

# Seed the freshly inserted row with the same values as the row directly
# below it (which now holds what used to be row 78), then overwrite the
# date with the new week's value.
$srcRow = $newRow + 1
$ws.Range("A" + $srcRow + ":T" + $srcRow).Copy()
$ws.Range("A" + $newRow + ":T" + $newRow).PasteSpecial()

$ws.Range("D" + $newRow).Value = 44932
